$d = $word.ActiveDocument

# --- 1. "CDN: ..." paragraph becomes the new "Difference between a framework
#        and a library?" question, followed by two brand-new paragraphs: the
#        library-vs-framework explanation and the (lightly reworded) CDN blurb.

$p1 = $d.Paragraphs(1)
$p1.Range.Text = "Difference between a framework and a library?"

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "A library is a piece of code that can be called to perform a specific task. Whereas a framework provides a basic skeleton or structure to your application. A framework contains libraries. A library can be applied to modify a single element in the application without having to disturb others."

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(3)
$p3.Range.Text = "CDN: Content delivery networks help to import the react code from the server and make our project configured to use react."

# --- 2. "For example in JS ..." gains a comma: "For example, in JS ..."

$find = $d.Content.Find
$find.Execute("For example in JS there's a creteElement API that lets you create an html element using JS.", $true, $false, $false, $false, $false, $true, 1, $false, "For example, in JS there's a creteElement API that lets you create an html element using JS.", 2) | Out-Null
